# إضافة حدث جديد في Card11 by admin at 2025-12-08 12:07:55
#
# This script:
#   1. Fills the previously-blank cells in row 18 of the "Card11" sheet
#      (B,C,D,E,F,G,H,I,J,K,M,P) with the text "nan" (matching how the
#      source data-export tool stamps empty numeric/text fields).
#   2. Appends a new event row (row 19) recording a new service entry:
#        Date (L19)    : 8\3\2025
#        Correction(N19): تم سن الفلاتس وتغيير الجرائد الخلفيه (1_5_8) ومعايره
#        Serviced by(O19): الخبير
#      with card number "11" in column A, and all other columns left
#      genuinely empty (no value at all).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card11")

# --- 1. Row 18: stamp the still-empty columns with the literal text "nan" ---
$row18Cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "M", "P")
foreach ($col in $row18Cols) {
    $ws.Range($col + "18").Value = "nan"
}

# --- 2. Row 19: new service-log entry ---
# Column A holds the card number "11". It is purely numeric text, so force
# it to be stored as text (matching the source column's inline-string type)
# without leaving any residual number-format/style change behind.
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "11"
$ws.Range("A19").ClearFormats()

# The remaining "unused" columns for this new row stay genuinely empty,
# but the source file still materialises an (empty) cell record for each
# of them. Touch-then-clear each one so a blank, default-style cell is
# written out instead of the column being omitted entirely.
$row19EmptyCols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "M", "P")
foreach ($col in $row19EmptyCols) {
    $cell = $ws.Range($col + "19")
    $cell.NumberFormat = "General"
    $cell.ClearFormats()
}

$ws.Range("L19").Value = "8\3\2025"
$ws.Range("N19").Value = "تم سن الفلاتس وتغيير الجرائد الخلفيه (1_5_8) ومعايره"
$ws.Range("O19").Value = "الخبير"
